$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$line1 = 'install.packages("C:\\Users\\Tom\\Documents\\Tom_Levers_Git_Repository\\R\\TomLeversRPackage", repos = NULL, type="source")'
$line2 = 'install.packages("tidyverse", repos = "http://cran.us.r-project.org")'

$ws.Range("B11").Value = $line1 + "`n" + $line2
